$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header text in A1 ("Nombre de logs" -> stored as a shared string)
$ws.Range("A1").Value = "Nombre de logs"

# Widen column A so the header text is fully visible
$ws.Columns.Item(1).ColumnWidth = 16.140625

# Leave the sheet with D6 as the active/selected cell
[void]$ws.Range("D6").Select()
